$d = $word.ActiveDocument

# --- Step 1: remove paragraphs that were dropped from the "what's next" list ---
# (delete bottom-up so indices stay valid)
#  50 Work on activity to activity navigation
#  52 Local stories
#  53 Set up grid view to save new stories
#  54 Redesign create a story activity
#  55 Create fragment for My stories, local stories, invited stories(done)
#  56 Remove auto rotate(done)
#  58 Pressing back in the main activity should close the app
#  59 Either make start a story page portrait only or add scrollview
$toRemove = @(59, 58, 56, 55, 54, 53, 52, 50)
foreach ($idx in $toRemove) {
    $d.Paragraphs.Item($idx).Range.Delete()
}

# --- Step 2: drop the old "_GoBack" bookmark that sat on the "Not Needed" paragraph;
#     it reappears below on the new "story repetition" line instead ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- Step 3: insert the new batch of "what's next" items (plus 2 extra blank lines)
#     right after the "(Activity lifecycle)" paragraph ---
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:lastRenderedPageBreak/><w:t>push notifications</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">story repetition </w:t></w:r><w:bookmarkStart w:id="1" w:name="_GoBack"/><w:bookmarkEnd w:id="1"/></w:p><w:p><w:r><w:t>reset forgotten password</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">apply different story </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>sizes(</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve">50 pages, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>etc</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>)</w:t></w:r></w:p><w:p><w:r><w:t>character type count for “what happens next” in story</w:t></w:r></w:p><w:p><w:r><w:t>story count for how many pages left</w:t></w:r></w:p><w:p><w:r><w:t>admin end story</w:t></w:r></w:p><w:p><w:r><w:t>real time dialog loading</w:t></w:r></w:p><w:p><w:r><w:t>block user from posting after post</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">on story </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>listview</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> click – show user info</w:t></w:r><w:r><w:t xml:space="preserve"> or highlight</w:t></w:r></w:p><w:p><w:r><w:t>clean up</w:t></w:r></w:p><w:p/><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$anchor = $d.Paragraphs.Item(52)
$insertAt = $d.Range($anchor.Range.End - 1, $anchor.Range.End - 1)
$insertAt.InsertXML($xml)
